{"js": "// The document contains a single 5x3 table of lattice-multiplication\n// exercises. Each cell holds one run with 5 lines (separated by <w:br/>,\n// represented here as a vertical-tab \"\\v\" in the Office.js text model):\n//   \"AB x CD\"\n//   \"  C    D\"\n//   \"  ----\"\n//   \"A|    |\"\n//   \"B|    |\"\n// This script replaces the 15 exercises with a new set of numbers,\n// keeping the table shape / formatting (run properties, cell widths, \u2026)\n// untouched.\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"values\");\nawait context.sync();\n\nconst rowCount = table.values.length;\nconst colCount = table.values[0].length;\n\nconst after = [\n  [\n    \"66 x 56\\v  5    6\\v  ----\\v6|    |\\v6|    |\",\n    \"98 x 92\\v  9    2\\v  ----\\v9|    |\\v8|    |\",\n    \"14 x 72\\v  7    2\\v  ----\\v1|    |\\v4|    |\",\n  ],\n  [\n    \"55 x 99\\v  9    9\\v  ----\\v5|    |\\v5|    |\",\n    \"17 x 78\\v  7    8\\v  ----\\v1|    |\\v7|    |\",\n    \"33 x 29\\v  2    9\\v  ----\\v3|    |\\v3|    |\",\n  ],\n  [\n    \"46 x 49\\v  4    9\\v  ----\\v4|    |\\v6|    |\",\n    \"77 x 67\\v  6    7\\v  ----\\v7|    |\\v7|    |\",\n    \"10 x 49\\v  4    9\\v  ----\\v1|    |\\v0|    |\",\n  ],\n  [\n    \"32 x 48\\v  4    8\\v  ----\\v3|    |\\v2|    |\",\n    \"13 x 47\\v  4    7\\v  ----\\v1|    |\\v3|    |\",\n    \"25 x 26\\v  2    6\\v  ----\\v2|    |\\v5|    |\",\n  ],\n  [\n    \"95 x 38\\v  3    8\\v  ----\\v9|    |\\v5|    |\",\n    \"32 x 70\\v  7    0\\v  ----\\v3|    |\\v2|    |\",\n    \"15 x 79\\v  7    9\\v  ----\\v1|    |\\v5|    |\",\n  ],\n];\n\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    const range = cell.body.paragraphs.getFirst().getRange();\n    range.insertText(after[r][c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single 5x3 table of lattice-multiplication\n# exercises. Each cell's Range.Text is 5 lines joined by the Word\n# line-break character (Chr(11), vertical tab -> <w:br/> on save):\n#   \"AB x CD\"\n#   \"  C    D\"\n#   \"  ----\"\n#   \"A|    |\"\n#   \"B|    |\"\n# This script replaces the 15 exercises with a new set of numbers,\n# keeping the table shape / formatting untouched.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$cellTexts = @(\n    ('66 x 56' + [char]11 + '  5    6' + [char]11 + '  ----' + [char]11 + '6|    |' + [char]11 + '6|    |'),\n    ('98 x 92' + [char]11 + '  9    2' + [char]11 + '  ----' + [char]11 + '9|    |' + [char]11 + '8|    |'),\n    ('14 x 72' + [char]11 + '  7    2' + [char]11 + '  ----' + [char]11 + '1|    |' + [char]11 + '4|    |'),\n    ('55 x 99' + [char]11 + '  9    9' + [char]11 + '  ----' + [char]11 + '5|    |' + [char]11 + '5|    |'),\n    ('17 x 78' + [char]11 + '  7    8' + [char]11 + '  ----' + [char]11 + '1|    |' + [char]11 + '7|    |'),\n    ('33 x 29' + [char]11 + '  2    9' + [char]11 + '  ----' + [char]11 + '3|    |' + [char]11 + '3|    |'),\n    ('46 x 49' + [char]11 + '  4    9' + [char]11 + '  ----' + [char]11 + '4|    |' + [char]11 + '6|    |'),\n    ('77 x 67' + [char]11 + '  6    7' + [char]11 + '  ----' + [char]11 + '7|    |' + [char]11 + '7|    |'),\n    ('10 x 49' + [char]11 + '  4    9' + [char]11 + '  ----' + [char]11 + '1|    |' + [char]11 + '0|    |'),\n    ('32 x 48' + [char]11 + '  4    8' + [char]11 + '  ----' + [char]11 + '3|    |' + [char]11 + '2|    |'),\n    ('13 x 47' + [char]11 + '  4    7' + [char]11 + '  ----' + [char]11 + '1|    |' + [char]11 + '3|    |'),\n    ('25 x 26' + [char]11 + '  2    6' + [char]11 + '  ----' + [char]11 + '2|    |' + [char]11 + '5|    |'),\n    ('95 x 38' + [char]11 + '  3    8' + [char]11 + '  ----' + [char]11 + '9|    |' + [char]11 + '5|    |'),\n    ('32 x 70' + [char]11 + '  7    0' + [char]11 + '  ----' + [char]11 + '3|    |' + [char]11 + '2|    |'),\n    ('15 x 79' + [char]11 + '  7    9' + [char]11 + '  ----' + [char]11 + '1|    |' + [char]11 + '5|    |')\n)\n\n$rows = $tbl.Rows.Count\n$cols = $tbl.Columns.Count\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $cellTexts[$i]\n        $i++\n    }\n}\n"}
